# Apply the "old runs" documentation update to HW4/StatisticalAnalysis.xlsx
# (Data_Acc_Kernels sheet): fill in the Action 1.2 / Action 1.3 result
# tables, refresh the Action 1.4 table with newer numbers, and drop the
# now-unused "16/8" label cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data_Acc_Kernels")

# --- Remove the stray "16/8" note next to the old Action 1.4 table ---
$ws.Range("I12").ClearContents()

# --- Build header rows (row 11) for the two new tables by copying the
#     formatting already used for the existing G11:H11 header ---
$ws.Range("G11:H11").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("D11").PasteSpecial(-4122)

# --- Copy the row formatting (labels s=2, numbers s=3, average s=6/7)
#     down into the two new tables from the existing Action 1.4 table ---
$ws.Range("G12:H16").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("D12").PasteSpecial(-4122)

$ws.Range("G17:H17").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("D17").PasteSpecial(-4122)

# --- Table 1: Action 1.2 (columns A:B) ---
$ws.Range("A11").Value = "Laplace_2d_OpenAcc"
$ws.Range("B11").Value = "Run Time (s)"

$ws.Range("A12").Value = "Run 1"
$ws.Range("B12").Value = 209.35736700000001
$ws.Range("A13").Value = "Run 2"
$ws.Range("B13").Value = 201.06361200000001
$ws.Range("A14").Value = "Run 3"
$ws.Range("B14").Value = 215.540524
$ws.Range("A15").Value = "Run 4"
$ws.Range("B15").Value = 214.03624300000001
$ws.Range("A16").Value = "Run 5"
$ws.Range("B16").Value = 215.44611499999999

$ws.Range("A17").Value = "Average (s)"
$ws.Range("B17").Formula = "=AVERAGE(B12:B16)"

# --- Table 2: Action 1.3 (columns D:E) ---
$ws.Range("D11").Value = "Laplace_2d_OpenAcc"
$ws.Range("E11").Value = "Run Time (s)"

$ws.Range("D12").Value = "Run 1"
$ws.Range("E12").Value = 10.51169
$ws.Range("D13").Value = "Run 2"
$ws.Range("E13").Value = 9.6477319999999995
$ws.Range("D14").Value = "Run 3"
$ws.Range("E14").Value = 10.839015
$ws.Range("D15").Value = "Run 4"
$ws.Range("E15").Value = 10.424659999999999
$ws.Range("D16").Value = "Run 5"
$ws.Range("E16").Value = 9.5849729999999997

$ws.Range("D17").Value = "Average (s)"
$ws.Range("E17").Formula = "=AVERAGE(E12:E16)"

# --- Table 3: Action 1.4 (columns G:H) refreshed with new measurements ---
$ws.Range("H12").Value = 10.170871
$ws.Range("H13").Value = 9.3651929999999997
$ws.Range("H14").Value = 10.568781
$ws.Range("H15").Value = 10.236471
$ws.Range("H16").Value = 9.4068769999999997
# H17 keeps its existing =AVERAGE(H12:H16) formula; it recalculates itself.

# --- Update the view so the Action 1.4 table area is in focus/selected ---
$ws.Range("G11:H17").Select()
